$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidated values: each token's name/type-line/abilities/power-toughness
# rows are collapsed into a single row formatted as a Python-style tuple string.
$values = @(
    "('Beast', ['Token Creature — Beast', '8/8'])",
    "('Dragon', ['Token Creature — Dragon', 'Flying', '4/4'])",
    "('Goblin', ['Token Creature — Goblin', '1/1'])",
    "('Homunculus', ['Token Artifact Creature — Homunculus', '0/1'])",
    "('Ooze', ['Token Creature — Ooze', '*/*'])",
    "('Saproling', ['Token Creature — Saproling', '1/1'])",
    "('Skeleton', ['Token Creature — Skeleton', '{B}: Regenerate this creature.', '1/1'])",
    "('Soldier', ['Token Creature — Soldier', '1/1'])",
    "('Thopter', ['Token Artifact Creature — Thopter', 'Flying', '1/1'])",
    "('Zombie', ['Token Creature — Zombie', '2/2'])"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# The old data occupied rows 2-34; after consolidating into rows 2-11, remove
# the now-stale rows 12-34 and shift the (now empty) cells below up.
$xlShiftUp = -4162
$ws.Range("A12:A34").Delete($xlShiftUp) | Out-Null
